{"js": "// \"nos permitir\u00e1\" -> \"va a permitir\" inside the paragraph that begins\n// \"En resumen este ejemplo completo ...\". The original run (\"ste ejemplo\n// completo nos permitir\u00e1 repasar ...\") gets split into three runs around\n// the replaced phrase, all sharing the same run formatting.\nconst body = context.document.body;\nconst results = body.search(\"nos permitir\u00e1\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Text \"nos permitir\u00e1\" not found.');\n}\n\nconst target = results.items[0];\ntarget.insertText(\"va a permitir\", Word.InsertLocation.replace);\n// Toggling a character property and back forces Word to keep the newly\n// inserted text as its own run (matching the split the real edit produced)\n// instead of silently re-merging it with its identically formatted\n// neighbours.\ntarget.font.bold = true;\nawait context.sync();\ntarget.font.bold = false;\nawait context.sync();\n", "ps1": "# \"nos permitir\u00e1\" -> \"va a permitir\" inside the paragraph that begins\n# \"En resumen este ejemplo completo ...\". The original run (\"ste ejemplo\n# completo nos permitir\u00e1 repasar ...\") gets split into three runs around\n# the replaced phrase, all sharing the same run formatting.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"nos permitir\u00e1\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 0\n\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Text 'nos permitir\u00e1' not found\"\n}\n\n$rng = $find.Parent\n$rng.Text = \"va a permitir\"\n# Toggling a character property and back forces Word to keep the newly\n# inserted text as its own run (matching the split the real edit produced)\n# instead of silently re-merging it with its identically formatted\n# neighbours.\n$rng.Bold = 1\n$rng.Bold = 0\n"}
